$wb = $excel.ActiveWorkbook

# zh-cn sheet: update Correspond Handoff/Handback DateTime for row 2
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-11 08:44:30"
$wsZhCn.Range("H2").Value = "2016-03-11 08:44:47"

# de-de sheet: update Correspond Handoff/Handback DateTime for row 2
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-11 08:44:33"
$wsDeDe.Range("H2").Value = "2016-03-11 08:44:52"
